# Rewrite the simple `{ m:'doc.html'.fromHTMLURI() }` field (begin/instrText/end
# field-char run soup) into the same text expressed as plain `w:t` literal runs,
# keeping the `_GoBack` bookmark that sits in the middle of the field code.
#
# (Parser switched to TokenIteratorFieldRewriterSplit, which expects the M2Doc
# tokens to already be present as literal paragraph text instead of a Word field.)

$d = $word.ActiveDocument

# Locate the field whose code looks like the M2Doc "fromHTMLURI" query, and the
# paragraph that contains it (done by range containment so this does not rely
# on a hard-coded paragraph index).
$target = $null
foreach ($f in $d.Fields) {
    if ($f.Code.Text -match "fromHTMLURI") {
        $target = $f
        break
    }
}

$codeStart = $target.Code.Start
$codeEnd = $target.Code.End

$hostPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -le $codeStart -and $p.Range.End -ge $codeEnd) {
        $hostPara = $p
        break
    }
}

# Range covering the whole paragraph content (field-char begin .. field-char
# end) but excluding the trailing paragraph mark, so InsertXML replaces only
# the paragraph's runs and keeps the <w:p> itself (and its attributes) intact.
$start = $hostPara.Range.Start
$end = $hostPara.Range.End - 1
$fieldRange = $d.Range($start, $end)

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">
<w:r><w:t>{</w:t></w:r>
<w:r w:rsidR="00DE6D5A"><w:t>m</w:t></w:r>
<w:r w:rsidR="002033E1"><w:t>:</w:t></w:r>
<w:r w:rsidR="008B76C9"><w:t>'</w:t></w:r>
<w:r w:rsidR="00E806A4"><w:t>doc.html</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r w:rsidR="008B76C9"><w:t>'.fromHTMLURI()</w:t></w:r>
<w:r><w:t xml:space="preserve">}</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$fieldRange.InsertXML($xml)
